# Update cryptos price list: refresh Price (D) and Volume(1h) (E) columns
# for rows 2-51 with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.721.54"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "3.795.78"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "700.80"
$ws.Range("E5").Value = "  +5.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.36"
$ws.Range("E6").Value = "  +3.32%  "
$ws.Range("D7").Value = "3.795.11"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.23"
$ws.Range("E11").Value = "  +4.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000256"
$ws.Range("E13").Value = "  +6.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.08"
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("D15").Value = "4.435.60"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "3.784.36"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "70.630.71"
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.69"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.20"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.17"
$ws.Range("E21").Value = "  +15.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "479.15"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.711"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.93"
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000142"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.33"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.45"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.15"
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("D29").Value = "3.944.92"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.12"
$ws.Range("E31").Value = "  +12.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.55"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.28"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  +8.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.40"
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.23"
$ws.Range("E36").Value = "  +3.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.102"
$ws.Range("E38").Value = "  +1.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.39"
$ws.Range("E39").Value = "  +2.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.98"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.23"
$ws.Range("E41").Value = "  +9.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.977"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000319"
$ws.Range("E45").Value = "  +16.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.31"
$ws.Range("E46").Value = "  +3.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.86"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "44.54"
$ws.Range("E48").Value = "  -3.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.39"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.300"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.62"
$ws.Range("E51").Value = "  +1.73%  "
